# Mint_Alias.xlsx - add "Kazakhstan Mint" / Ust-Kamenogorsk row to the
# "Information" table (Таблица3), per meshok.net listing #71561403 (10.10.2024).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Information")

# The table currently spans A1:H40 (header + 39 data rows). Grow it by one
# row: copy the formatting of the last data row (row 40) down into row 41 so
# the new row picks up the same borders/alignment as its neighbours.
$ws.Range("A40:H40").Copy()
$ws.Range("A41:H41").PasteSpecial(-4122) | Out-Null

# Country column (D) on row 40 happens to carry the "last row" border
# variant; the freshly typed row instead matches the plain interior style
# used by column C, so pull that one cell's formatting in separately.
$ws.Range("C40").Copy()
$ws.Range("D41").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New record: ID 40 - Kazakhstan Mint (Ust-Kamenogorsk), founded 1992.
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = "Kazakhstan Mint"
$ws.Range("C41").Value = "Қазақстан теңге сарайы"
$ws.Range("D41").Value = "Kazakhstan"
$ws.Range("E41").Value = "Ust-Kamenogorsk"
$ws.Range("F41").Value = "QUB"
$ws.Range("G41").Value = 1992

# Extend the ListObject (table) + its AutoFilter to include the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:H41")) | Out-Null

# Match the author's final selection state on the sheet.
[void]$ws.Range("J44").Select()
